$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "pontos notáveis - incremento na tabela de ranking"
# Columns E and F (rows 2-7) store percentage-formatted values as raw
# fractions (e.g. 0.79...). Rescale them by a factor of 100 so they line
# up with the rest of the ranking increments.
$rng = $ws.Range("E2:F7")
for ($i = 1; $i -le $rng.Rows.Count; $i++) {
    for ($j = 1; $j -le $rng.Columns.Count; $j++) {
        $cell = $rng.Cells.Item($i, $j)
        $cell.Value2 = $cell.Value2 * 100
    }
}
